$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 57 and 58 had their match data (columns F..V) swapped.
#    Columns A..E (index / pais / torneio / temporada / data_partida) stay
#    put since both matches were played on the same date; only the
#    match-specific data (teams, goals, odds, timestamps, url) moves.
# ---------------------------------------------------------------------------
$cols = 6..22
$row57 = @{}
$row58 = @{}
foreach ($c in $cols) {
    $row57[$c] = $ws.Cells.Item(57, $c).Value2
    $row58[$c] = $ws.Cells.Item(58, $c).Value2
}
foreach ($c in $cols) {
    $ws.Cells.Item(57, $c).Value = $row58[$c]
    $ws.Cells.Item(58, $c).Value = $row57[$c]
}

# ---------------------------------------------------------------------------
# 2) Two new match rows were appended at the bottom (110 and 111).
#    Clone the formatting of the last existing data row (109) so the new
#    rows reuse the same styles (bold/bordered index column, date column
#    number format) instead of creating new style entries.
# ---------------------------------------------------------------------------
$ws.Range("A109:E109").Copy()
$ws.Range("A110:E111").PasteSpecial(-4122)

# Row 110: Frosinone 2 - 1 Empoli
$ws.Cells.Item(110, 1).Value = 109
$ws.Cells.Item(110, 2).Value = "italy"
$ws.Cells.Item(110, 3).Value = "serie-a"
$ws.Cells.Item(110, 4).Value = "2023-2024"
$ws.Cells.Item(110, 5).Value = 45236.77083333334
$ws.Cells.Item(110, 6).Value = "Frosinone"
$ws.Cells.Item(110, 7).Value = 2
$ws.Cells.Item(110, 8).Value = "Empoli"
$ws.Cells.Item(110, 9).Value = 1
$ws.Cells.Item(110, 10).Value = 1.98
$ws.Cells.Item(110, 11).Value = "23/10/2023 15:49"
$ws.Cells.Item(110, 12).Value = 2.18
$ws.Cells.Item(110, 13).Value = "06/11/2023 18:04"
$ws.Cells.Item(110, 14).Value = 3.56
$ws.Cells.Item(110, 15).Value = "23/10/2023 15:49"
$ws.Cells.Item(110, 16).Value = 3.59
$ws.Cells.Item(110, 17).Value = "06/11/2023 18:29"
$ws.Cells.Item(110, 18).Value = 3.68
$ws.Cells.Item(110, 19).Value = "23/10/2023 15:49"
$ws.Cells.Item(110, 20).Value = 3.44
$ws.Cells.Item(110, 21).Value = "06/11/2023 18:29"
$ws.Cells.Item(110, 22).Value = "https://www.betexplorer.com/football/italy/serie-a/frosinone-empoli/fkMSidlb/"

# Row 111: Torino 2 - 1 Sassuolo
$ws.Cells.Item(111, 1).Value = 110
$ws.Cells.Item(111, 2).Value = "italy"
$ws.Cells.Item(111, 3).Value = "serie-a"
$ws.Cells.Item(111, 4).Value = "2023-2024"
$ws.Cells.Item(111, 5).Value = 45236.86458333334
$ws.Cells.Item(111, 6).Value = "Torino"
$ws.Cells.Item(111, 7).Value = 2
$ws.Cells.Item(111, 8).Value = "Sassuolo"
$ws.Cells.Item(111, 9).Value = 1
$ws.Cells.Item(111, 10).Value = 1.89
$ws.Cells.Item(111, 11).Value = "22/10/2023 12:02"
$ws.Cells.Item(111, 12).Value = 2.09
$ws.Cells.Item(111, 13).Value = "06/11/2023 20:41"
$ws.Cells.Item(111, 14).Value = 3.69
$ws.Cells.Item(111, 15).Value = "22/10/2023 12:02"
$ws.Cells.Item(111, 16).Value = 3.45
$ws.Cells.Item(111, 17).Value = "06/11/2023 20:37"
$ws.Cells.Item(111, 18).Value = 4.23
$ws.Cells.Item(111, 19).Value = "22/10/2023 12:02"
$ws.Cells.Item(111, 20).Value = 3.86
$ws.Cells.Item(111, 21).Value = "06/11/2023 20:43"
$ws.Cells.Item(111, 22).Value = "https://www.betexplorer.com/football/italy/serie-a/torino-sassuolo/OUDqlEmN/"
